$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = "30.191.78"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(2, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.91%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = "1.850.20"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(3, 5)
$c.NumberFormat = "@"
$c.Value = "  -2.34%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value = "0.9995"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "236.21"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(5, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.81%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "0.9991"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(6, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.11%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(7, 5)
$c.NumberFormat = "@"
$c.Value = "  -2.51%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = "0.2816"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(8, 5)
$c.NumberFormat = "@"
$c.Value = "  -4.04%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(9, 5)
$c.NumberFormat = "@"
$c.Value = "  -3.25%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = "1.856.34"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(10, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.90%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "0.07303"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(11, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.66%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = "16.34"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(12, 5)
$c.NumberFormat = "@"
$c.Value = "  -4.42%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = "5.133"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(13, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.07%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "87.26"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(14, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.98%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "0.6450"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(15, 5)
$c.NumberFormat = "@"
$c.Value = "  -3.20%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "30.137.32"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = "13.23"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(17, 5)
$c.NumberFormat = "@"
$c.Value = "  -1.90%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = "0.9996"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(18, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.11%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "0.000007626"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(19, 5)
$c.NumberFormat = "@"
$c.Value = "  -2.54%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(20, 2)
$c.NumberFormat = "@"
$c.Value = "BitcoinCash"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(20, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = "226.19"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(20, 5)
$c.NumberFormat = "@"
$c.Value = "  +18.96%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(21, 2)
$c.NumberFormat = "@"
$c.Value = "WrappedliquidstakedEther2.0"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(21, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "2.093.38"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(21, 5)
$c.NumberFormat = "@"
$c.Value = "  -2.41%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(22, 2)
$c.NumberFormat = "@"
$c.Value = "BinanceUSD"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(22, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "1.000"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(22, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.06%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(23, 2)
$c.NumberFormat = "@"
$c.Value = "Uniswap"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(23, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "5.272"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(23, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.79%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "6.095"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(24, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.89%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "9.212"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(25, 5)
$c.NumberFormat = "@"
$c.Value = "  -2.98%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = "163.85"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(26, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.28%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "18.49"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(27, 5)
$c.NumberFormat = "@"
$c.Value = "  +1.02%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = "1.915"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(28, 5)
$c.NumberFormat = "@"
$c.Value = "  -1.06%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "1.428"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(29, 5)
$c.NumberFormat = "@"
$c.Value = "  -2.40%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(30, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.22%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(31, 5)
$c.NumberFormat = "@"
$c.Value = "  -2.65%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "3.959"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(32, 5)
$c.NumberFormat = "@"
$c.Value = "  -3.00%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = "0.05021"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(33, 5)
$c.NumberFormat = "@"
$c.Value = "  -3.78%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "0.7407"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(34, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.32%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = "1.137"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(35, 5)
$c.NumberFormat = "@"
$c.Value = "  +3.36%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(36, 5)
$c.NumberFormat = "@"
$c.Value = "  -1.18%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = "0.01824"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(37, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.19%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = "2.607"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(38, 5)
$c.NumberFormat = "@"
$c.Value = "  -2.57%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(39, 2)
$c.NumberFormat = "@"
$c.Value = "TrustWalletToken"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(39, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "0.9023"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(39, 5)
$c.NumberFormat = "@"
$c.Value = "  -2.25%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(40, 2)
$c.NumberFormat = "@"
$c.Value = "RenderToken"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(40, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "2.055"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(40, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.32%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(41, 2)
$c.NumberFormat = "@"
$c.Value = "Quant"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(41, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "106.89"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(41, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.48%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(42, 2)
$c.NumberFormat = "@"
$c.Value = "FraxShare"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(42, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "5.915"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(42, 5)
$c.NumberFormat = "@"
$c.Value = "  -0.75%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = "0.4246"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(43, 5)
$c.NumberFormat = "@"
$c.Value = "  -3.69%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "0.9985"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(44, 5)
$c.NumberFormat = "@"
$c.Value = "  +0.47%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "7.425"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(45, 5)
$c.NumberFormat = "@"
$c.Value = "  -2.43%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "0.1311"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(46, 5)
$c.NumberFormat = "@"
$c.Value = "  -5.41%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "1.556"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(47, 5)
$c.NumberFormat = "@"
$c.Value = "  +9.98%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = "63.87"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(48, 5)
$c.NumberFormat = "@"
$c.Value = "  -6.37%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(49, 2)
$c.NumberFormat = "@"
$c.Value = "Elrond"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(49, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = "34.33"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(49, 5)
$c.NumberFormat = "@"
$c.Value = "  -1.82%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(50, 2)
$c.NumberFormat = "@"
$c.Value = "EnergySwap"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(50, 3)
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = "8.766"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(50, 5)
$c.NumberFormat = "@"
$c.Value = "  -2.88%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = "0.05675"
$c.NumberFormat = "General"
$c.Style = "Normal"

$c = $ws.Cells.Item(51, 5)
$c.NumberFormat = "@"
$c.Value = "  -2.53%  "
$c.NumberFormat = "General"
$c.Style = "Normal"

